$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.800.86"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "2.583.92"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.87"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.90"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.352"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.97"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "3.047.24"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "62.669.08"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "2.588.69"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.24"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.93"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.35"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.27"
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("D24").Value = "2.706.62"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("E25").Value = "  -1.99%  "
$ws.Range("E26").Value = "  -2.43%  "
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.94"
$ws.Range("E29").Value = "  +3.16%  "
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.93"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "458.75"
$ws.Range("E33").Value = "  +8.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "176.78"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("E35").Value = "  +3.69%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.46"
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "157.74"
$ws.Range("E42").Value = "  +4.62%  "
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.17"
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("E45").Value = "  +4.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0537"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.42"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("E51").Value = "  -1.11%  "
